$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "N" column (year 2020) -------------------------------------------

# N3: blank cell on the thin divider row, same look as the rest of row 3
# (Times New Roman 9pt body font + medium black bottom border)
$ws.Range("N3").Font.Name = "Times New Roman"
$ws.Range("N3").Font.Size = 9
$b3 = $ws.Range("N3").Borders.Item(9)
$b3.LineStyle = 1
$b3.Weight = -4138
$b3.Color = 0

# N4: header cell "2020" - bold 9pt + medium black bottom border
$ws.Range("N4").Value = 2020
$ws.Range("N4").Font.Name = "Times New Roman"
$ws.Range("N4").Font.Size = 9
$ws.Range("N4").Font.Bold = $true
$b4 = $ws.Range("N4").Borders.Item(9)
$b4.LineStyle = 1
$b4.Weight = -4138
$b4.Color = 0

# N5: blank cell under the section header row - regular 9pt, no border
$ws.Range("N5").Font.Name = "Times New Roman"
$ws.Range("N5").Font.Size = 9

# N6:N9 data values - regular 9pt, no border
$ws.Range("N6").Value = 1713
$ws.Range("N7").Value = 1
$ws.Range("N8").Value = 379
$ws.Range("N9").Value = 180
$ws.Range("N6:N9").Font.Name = "Times New Roman"
$ws.Range("N6:N9").Font.Size = 9

# N10: last data value - regular 9pt + medium black bottom border (table close)
$ws.Range("N10").Value = 798
$ws.Range("N10").Font.Name = "Times New Roman"
$ws.Range("N10").Font.Size = 9
$b10 = $ws.Range("N10").Borders.Item(9)
$b10.LineStyle = 1
$b10.Weight = -4138
$b10.Color = 0

# --- Restore the saved cursor/selection state ------------------------------
$ws.Range("L22").Select()
